$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.914.37"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.752.34"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.61%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.58%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.64%  "

# Row 9
$ws.Range("E9").Value = "  -3.59%  "

# Row 10
$ws.Range("E10").Value = "  +1.02%  "

# Row 11
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.64"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -16.61%  "

# Row 12
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.382"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.23%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.238.69"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.50%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.50"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.11%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.529.93"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.86%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000150"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.20%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.756.03"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.32%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.10"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.29%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.06%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.34%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.72"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.40%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.534"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.60%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.03"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.34%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.85%  "

# Row 26
$ws.Range("E26").Value = "  +0.48%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.88%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0906"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.97%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.94"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.92%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.06"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.88%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.25"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.45%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.93"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.62%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.16"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.77%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.44%  "

# Row 35
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.45"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.37%  "

# Row 37
$ws.Range("E37").Value = "  -1.74%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.985"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.91%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.17"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.83%  "

# Row 40
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "331.72"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.60%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.14"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.36%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.98"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.95%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.70%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0586"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.42%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.49"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.66%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0254"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.44%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.08"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.625"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.04%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.101"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.08%  "

# Row 50
$ws.Range("E50").Value = "  -0.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.02"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.04%  "
